# Automated daily update of the "剩余" (remaining days) and "开始时间" (start date)
# columns on Sheet1. For every data row (2..99):
#   - if the remaining-days value (column E) is 1, the cycle is renewed:
#       E -> 10  and  F (start date, stored as an 8-digit yyyymmdd number) -> F + 10
#   - otherwise the remaining-days value is simply decremented by 1 (E -> E - 1)
# Row 36 is skipped because its start-date value is malformed (9 digits instead of
# the expected 8), matching the source data that was left untouched by this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = 99
$skipRows = @(36)

for ($row = 2; $row -le $lastRow; $row++) {
    if ($skipRows -contains $row) {
        continue
    }

    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq 1) {
        $eCell.Value2 = 10
        $fCell.Value2 = $fVal + 10
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
